# Insert a new data row at row 160 (pushes existing rows 160-302 down to 161-303)
# and populate it with the new weekly price record for
# "Vega Modelo de Temuco - Acelga".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(160).Insert()

$ws.Cells.Item(160, 1).Value = 10
$ws.Cells.Item(160, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(160, 3).Value = "La Araucanía"
$ws.Cells.Item(160, 4).Value = 44658
$ws.Cells.Item(160, 5).Value = 9
$ws.Cells.Item(160, 6).Value = 100112009
$ws.Cells.Item(160, 7).Value = "Acelga"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 65
$ws.Cells.Item(160, 11).Value = 9000
$ws.Cells.Item(160, 12).Value = 9000
$ws.Cells.Item(160, 13).Value = 9000
$ws.Cells.Item(160, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(160, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(160, 16).Value = 750
$ws.Cells.Item(160, 17).Value = 12
$ws.Cells.Item(160, 18).Value = "Hortaliza"
